$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D cells whose new values look numeric need to be pinned to Text format
# first, otherwise Excel auto-converts them to actual numbers (losing exact
# formatting / trailing zeros) instead of keeping them as literal strings.
$textCells = @("D5", "D6", "D9", "D11", "D13", "D16", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D34", "D36", "D38", "D39", "D40", "D41", "D44", "D45", "D47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.893.37"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "3.411.41"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "404.40"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "132.99"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("E7").Value = "  -2.88%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "0.687"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("E10").Value = "  -5.57%  "
$ws.Range("D11").Value = "41.87"
$ws.Range("E11").Value = "  -3.99%  "
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").Value = "8.40"
$ws.Range("E13").Value = "  -4.98%  "
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").Value = "3.388.12"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").Value = "11.64"
$ws.Range("E16").Value = "  +7.31%  "
$ws.Range("D17").Value = "61.876.64"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("E18").Value = "  -3.64%  "
$ws.Range("D19").Value = "0.0000139"
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("E20").Value = "  -5.44%  "
$ws.Range("D21").Value = "83.24"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "311.71"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "12.73"
$ws.Range("E23").Value = "  -3.48%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").Value = "4.82"
$ws.Range("E25").Value = "  +10.34%  "
$ws.Range("D26").Value = "29.60"
$ws.Range("E26").Value = "  -3.34%  "
$ws.Range("D27").Value = "7.94"
$ws.Range("E27").Value = "  +4.90%  "
$ws.Range("D28").Value = "8.06"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").Value = "2.75"
$ws.Range("E29").Value = "  +5.62%  "
$ws.Range("D30").Value = "0.172"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("D32").Value = "42.50"
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "11.35"
$ws.Range("E34").Value = "  -4.67%  "
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("D36").Value = "51.38"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "2.97"
$ws.Range("E38").Value = "  -2.56%  "
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  -6.84%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "0.322"
$ws.Range("E40").Value = "  +11.50%  "
$ws.Range("D41").Value = "139.74"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").Value = "3.94"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Value = "16.59"
$ws.Range("E45").Value = "  -5.49%  "
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "21.04"
$ws.Range("E47").Value = "  -5.27%  "
$ws.Range("D48").Value = "2.107.13"
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("E49").Value = "  +23.57%  "
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("E51").Value = "  +3.33%  "
